# Player.xlsx edit:
# "unify the conception of DataNode, DataTable, Entity."
#
# - Rename Property1 -> DataNode_1
# - Rename Property2 -> DataNode_2
# - Rename Record_Hero -> DataTable_Hero
# - Rename Record_Bag -> DataTable_Bag
# - Rename Record_CommPropertyValue -> DataTable_CommPropertyValue
# - Rename Record_Task -> DataTable_Task
# - Remove the Record_Building sheet entirely (and its now-unused
#   Building* shared strings / comments / vml drawing go with it)
# - Component sheet keeps its name
# - DataTable_Hero (formerly Record_Hero) becomes the active/selected sheet

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Property1").Name = "DataNode_1"
$wb.Worksheets.Item("Property2").Name = "DataNode_2"
$wb.Worksheets.Item("Record_Hero").Name = "DataTable_Hero"
$wb.Worksheets.Item("Record_Bag").Name = "DataTable_Bag"
$wb.Worksheets.Item("Record_CommPropertyValue").Name = "DataTable_CommPropertyValue"
$wb.Worksheets.Item("Record_Task").Name = "DataTable_Task"

# Drop the old building/record sheet - unifying the naming scheme means
# it no longer belongs alongside the DataNode/DataTable sheets.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Record_Building").Delete() | Out-Null
$excel.DisplayAlerts = $true

# Make the hero data table the active tab, as in the target workbook.
$wb.Worksheets.Item("DataTable_Hero").Activate()
